$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.440.61'
$ws.Range("E2").Value = '  +1.06%  '
$ws.Range("D3").Value = '1.944.07'
$ws.Range("E3").Value = '  -0.96%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '243.36'
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("D6").Value = '0.615'
$ws.Range("E6").Value = '  -1.82%  '
$ws.Range("D7").Value = '58.14'
$ws.Range("E7").Value = '  -7.08%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '0.365'
$ws.Range("E9").Value = '  -1.85%  '
$ws.Range("D10").Value = '55.65'
$ws.Range("E10").Value = '  -0.74%  '
$ws.Range("D11").Value = '0.0833'
$ws.Range("E11").Value = '  +3.42%  '
$ws.Range("D12").Value = '0.103'
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("D13").Value = '0.821'
$ws.Range("E13").Value = '  -4.47%  '
$ws.Range("D14").Value = '21.51'
$ws.Range("E14").Value = '  -3.03%  '
$ws.Range("D15").Value = '2.232.64'
$ws.Range("E15").Value = '  -0.51%  '
$ws.Range("D16").Value = '13.55'
$ws.Range("E16").Value = '  -3.67%  '
$ws.Range("D17").Value = '5.23'
$ws.Range("E17").Value = '  -3.68%  '
$ws.Range("D18").Value = '1.954.21'
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("D19").Value = '36.363.79'
$ws.Range("E19").Value = '  +1.15%  '
$ws.Range("D20").Value = '69.68'
$ws.Range("E20").Value = '  -1.98%  '
$ws.Range("D21").Value = '0.0₃0862'
$ws.Range("E21").Value = '  +0.83%  '
$ws.Range("D22").Value = '229.41'
$ws.Range("E22").Value = '  -3.31%  '
$ws.Range("D23").Value = '5.06'
$ws.Range("E23").Value = '  -2.80%  '
$ws.Range("E24").Value = '  -0.29%  '
$ws.Range("D25").Value = '2.43'
$ws.Range("E25").Value = '  -4.60%  '
$ws.Range("D26").Value = '2.28'
$ws.Range("E26").Value = '  -0.55%  '
$ws.Range("E27").Value = '  -6.31%  '
$ws.Range("D28").Value = '161.64'
$ws.Range("E28").Value = '  +1.42%  '
$ws.Range("D29").Value = '19.39'
$ws.Range("E29").Value = '  -2.21%  '
$ws.Range("D30").Value = '0.126'
$ws.Range("E30").Value = '  -3.61%  '
$ws.Range("D31").Value = '0.117'
$ws.Range("E31").Value = '  -1.80%  '
$ws.Range("E32").Value = '  +0.68%  '
$ws.Range("D33").Value = '4.66'
$ws.Range("E33").Value = '  -3.85%  '
$ws.Range("D34").Value = '0.0625'
$ws.Range("E34").Value = '  +0.60%  '
$ws.Range("D35").Value = '4.26'
$ws.Range("E35").Value = '  -3.33%  '
$ws.Range("E36").Value = '  -1.66%  '
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("E38").Value = '  -2.59%  '
$ws.Range("D39").Value = '2.13'
$ws.Range("E39").Value = '  -6.56%  '
$ws.Range("D40").Value = '3.02'
$ws.Range("E40").Value = '  -3.29%  '
$ws.Range("D41").Value = '0.0985'
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("D42").Value = '2.85'
$ws.Range("E42").Value = '  +1.79%  '
$ws.Range("D43").Value = '1.17'
$ws.Range("E43").Value = '  -4.28%  '
$ws.Range("E44").Value = '  -1.53%  '
$ws.Range("D45").Value = '16.05'
$ws.Range("E45").Value = '  -1.04%  '
$ws.Range("D46").Value = '1.348.12'
$ws.Range("E46").Value = '  +0.76%  '
$ws.Range("E47").Value = '  -5.30%  '
$ws.Range("D48").Value = '87.60'
$ws.Range("E48").Value = '  -5.31%  '
$ws.Range("D49").Value = '7.16'
$ws.Range("E49").Value = '  -5.24%  '
$ws.Range("E50").Value = '  +1.12%  '
$ws.Range("D51").Value = '45.45'
$ws.Range("E51").Value = '  +3.44%  '
